$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 7921
$ws1.Range("F4").Value = 85
$ws1.Range("F5").Value = 16669
$ws1.Range("F7").Value = 596
$ws1.Range("F8").Value = 656
$ws1.Range("F12").Value = 787
$ws1.Range("F14").Value = 88
$ws1.Range("F15").Value = 351
$ws1.Range("F17").Value = 305
$ws1.Range("F18").Value = 144
$ws1.Range("F21").Value = 1103
$ws1.Range("F23").Value = 657
$ws1.Range("F24").Value = 2244
$ws1.Range("F25").Value = 769
$ws1.Range("F27").Value = 565

# Sheet 3: 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 485

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 485
$ws4.Range("F3").Value = 7921
$ws4.Range("F5").Value = 85
$ws4.Range("F7").Value = 16672
$ws4.Range("F9").Value = 596
$ws4.Range("F10").Value = 656
$ws4.Range("F18").Value = 787
$ws4.Range("F20").Value = 88
$ws4.Range("F21").Value = 351
$ws4.Range("F27").Value = 305
$ws4.Range("F28").Value = 144
$ws4.Range("F31").Value = 1103
$ws4.Range("F33").Value = 657
$ws4.Range("F34").Value = 2244
$ws4.Range("F35").Value = 769
$ws4.Range("F37").Value = 565
